# Html failu drukāšana bet nepareizajiem unikodiem.
#
# Updates the "Sketch parametrs" labels in column A of Sheet1: the old
# lower-case English-ish labels are replaced with properly capitalized /
# translated labels, and three rows (Author, Chapter ID, Source ID) that
# previously had no label get one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Sheet"
$ws.Range("A3").Value = "Title (even)"
$ws.Range("A4").Value = "Title (odd)"
$ws.Range("A5").Value = "Author"
$ws.Range("A6").Value = "Carry-over"
$ws.Range("A7").Value = "Polish"
$ws.Range("A8").Value = "French"
$ws.Range("A9").Value = "Estonian"
$ws.Range("A10").Value = "Flemish"
$ws.Range("A12").Value = "Greek"
$ws.Range("A13").Value = "Italian"
$ws.Range("A15").Value = "Comment"
$ws.Range("A16").Value = "Latin"
$ws.Range("A18").Value = "Chapter ID"
$ws.Range("A1").Value = "Nosaukums Sketch-ā un tooltip-os"
$ws.Range("A20").Value = "Remark"
$ws.Range("A22").Value = "Aramaic"
$ws.Range("A23").Value = "English"
$ws.Range("A24").Value = "Parallel"
$ws.Range("A26").Value = "German"
$ws.Range("A30").Value = "Source ID"
